{"js": "// Update the two-digit division problems in the table to the new set\n// of values. Each original expression occurs exactly once in the\n// document, so searching the body for each exact string and replacing\n// its matched range's text is safe and keeps run formatting intact.\n\nconst replacements = [\n  [\"96\u00f74=\", \"70\u00f76=\"],\n  [\"66\u00f77=\", \"35\u00f74=\"],\n  [\"35\u00f76=\", \"15\u00f73=\"],\n  [\"89\u00f76=\", \"89\u00f79=\"],\n  [\"98\u00f78=\", \"44\u00f74=\"],\n  [\"46\u00f79=\", \"29\u00f72=\"],\n  [\"76\u00f77=\", \"65\u00f79=\"],\n  [\"99\u00f76=\", \"59\u00f73=\"],\n  [\"61\u00f72=\", \"13\u00f72=\"],\n  [\"86\u00f74=\", \"51\u00f72=\"],\n  [\"11\u00f72=\", \"96\u00f76=\"],\n  [\"62\u00f78=\", \"75\u00f78=\"],\n  [\"32\u00f72=\", \"24\u00f79=\"],\n  [\"32\u00f73=\", \"64\u00f74=\"],\n  [\"34\u00f72=\", \"73\u00f79=\"],\n  [\"73\u00f72=\", \"62\u00f73=\"],\n  [\"20\u00f73=\", \"36\u00f77=\"],\n  [\"27\u00f72=\", \"76\u00f73=\"],\n  [\"98\u00f75=\", \"90\u00f79=\"],\n  [\"77\u00f78=\", \"92\u00f77=\"],\n  [\"32\u00f79=\", \"20\u00f79=\"],\n  [\"80\u00f73=\", \"14\u00f72=\"],\n  [\"28\u00f73=\", \"15\u00f79=\"],\n  [\"21\u00f74=\", \"17\u00f79=\"],\n  [\"12\u00f73=\", \"68\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit division problems in the table to the new set\n# of values. Each original expression is unique in the document, so a\n# straightforward Find/Replace (whole match, no wildcards) for each pair\n# is safe and will not collide with any other cell.\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"96\u00f74=\" = \"70\u00f76=\"\n    \"66\u00f77=\" = \"35\u00f74=\"\n    \"35\u00f76=\" = \"15\u00f73=\"\n    \"89\u00f76=\" = \"89\u00f79=\"\n    \"98\u00f78=\" = \"44\u00f74=\"\n    \"46\u00f79=\" = \"29\u00f72=\"\n    \"76\u00f77=\" = \"65\u00f79=\"\n    \"99\u00f76=\" = \"59\u00f73=\"\n    \"61\u00f72=\" = \"13\u00f72=\"\n    \"86\u00f74=\" = \"51\u00f72=\"\n    \"11\u00f72=\" = \"96\u00f76=\"\n    \"62\u00f78=\" = \"75\u00f78=\"\n    \"32\u00f72=\" = \"24\u00f79=\"\n    \"32\u00f73=\" = \"64\u00f74=\"\n    \"34\u00f72=\" = \"73\u00f79=\"\n    \"73\u00f72=\" = \"62\u00f73=\"\n    \"20\u00f73=\" = \"36\u00f77=\"\n    \"27\u00f72=\" = \"76\u00f73=\"\n    \"98\u00f75=\" = \"90\u00f79=\"\n    \"77\u00f78=\" = \"92\u00f77=\"\n    \"32\u00f79=\" = \"20\u00f79=\"\n    \"80\u00f73=\" = \"14\u00f72=\"\n    \"28\u00f73=\" = \"15\u00f79=\"\n    \"21\u00f74=\" = \"17\u00f79=\"\n    \"12\u00f73=\" = \"68\u00f79=\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $find.Execute(\n        $old,      # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $new,      # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
